# updated main GSC export data
# The "Chart" sheet (rolling Video-Indexing export window) rolled forward
# by one day: the oldest day (2025-10-23) drops off the front, the three
# newest remaining days have not been backfilled with data yet, and the
# last row's Impressions value is now a real number instead of the
# "not available" placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the oldest day (2025-10-23) - every following row shifts up one.
$ws.Rows.Item(2).Delete()

# The three oldest days still in the window (2025-10-24 .. 2025-10-26)
# have no "No video indexed" / "Video indexed" counts yet.
$ws.Range("B2:C4").Value = ""

# The last row (2026-01-19) now reports a real Impressions value instead
# of the earlier "not available" placeholder.
$ws.Range("D89").Value = 0
